$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'245.24"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = "'24.08"
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.Value = "'5.212"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'0.05793"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'6.509"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.Value = "'3.128"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.Value = "'0.8571"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.Value = "'0.06971"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = "'0.02872"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'0.09375"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.Value = "'3.744"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = "'0.001508"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = "'0.04704"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.Value = "'0.0006011"
$cell.Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$cell = $ws.Range("D19")
$cell.Value = "'0.006288"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.Value = "'0.001235"
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.Value = "'0.004534"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = "'0.00006901"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.Value = "'3.500"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = "'0.03654"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = "'0.006301"
$cell.Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$cell = $ws.Range("D43")
$cell.Value = "'0.003401"
$cell.Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"
$cell = $ws.Range("D44")
$cell.Value = "'0.007948"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'0.00005260"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = "'0.3201"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = "'0.002351"
$cell.Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$cell = $ws.Range("D49")
$cell.Value = "'0.00002100"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'0.0002000"
$cell.Style = "Normal"
